$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column O data for year 2021, matching the style of the existing
# N column (year 2020) cells.
$ws.Range("N4:N5").Copy()
$ws.Range("O4:O5").PasteSpecial(-4122)

$ws.Range("O4").Value = 2021
$ws.Range("O5").Value = 515

# Reset the view: scroll back to show column A, and move the selection.
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("P12").Select()
